$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.760.07"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.533.13"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.20"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.21"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.751.88"
$ws.Range("E12").Value = "  -1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.522.50"
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.748.52"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.84"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0680"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.24"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.55"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0452"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  +2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.359.95"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  -2.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.951"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0164"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.799"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.42"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.665.78"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.11"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0502"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0969"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0941"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  -0.07%  "
